# Add a new "pre_first_results" column (E) to Sheet1.
# TRUE for rows recorded before the first election results came in
# (through 2020-11-03 18:00, serial 44138.75), FALSE from the next
# recorded row onward (2020-11-03 19:00, serial 44138.791666666664) to
# the end of the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("E1").Value = "pre_first_results"

# Data rows: 2-217 -> TRUE (pre first results), 218-267 -> FALSE (post)
$ws.Range("E2:E217").Value = $true
$ws.Range("E218:E267").Value = $false

# Update selection to match the newly edited range
$ws.Range("E218:E267").Select()
